# Update the "dSF" column (F) values for several rows in the Jordan Hicks
# 2022 save-data sheet. This reflects a repull/push of data along with a
# recalculation of the mean that shifted several dSF figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -2
    7  = -1
    13 = 2
    15 = 1
    19 = 0
    22 = -1
    23 = -4
    24 = 4
    29 = 1
    30 = -4
    32 = -2
    33 = 2
    35 = -1
    36 = -1
    39 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
